$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: seed new shared-string values in the exact order needed so the workbook's
# shared string table (xl/sharedStrings.xml) ends up matching the target ordering.
$ws.Range("M4").Value2 = "10^6 Hz"
$ws.Range("H4").Value2 = "MHz"
$ws.Range("I2").Value2 = "Frequency"
$ws.Range("G2").Value2 = "*10^3Hz"
$ws.Range("G4").Value2 = "*10^6 Hz"
$ws.Range("H2").Value2 = "kHz"
$ws.Range("G3").Value2 = "1000 Hz"
$ws.Range("G5").Value2 = "1000000 Hz"
$ws.Range("G6").Value2 = "*10^9 Hz"
$ws.Range("G7").Value2 = "1000000000 Hz"
$ws.Range("H6").Value2 = "GHz"
$ws.Range("G8").Value2 = "rad/second"
$ws.Range("H8").Value2 = "rad*s^{-1}"
$ws.Range("G9").Value2 = "Pascal"
$ws.Range("H9").Value2 = "Pa"
$ws.Range("I9").Value2 = "Pressure"
$ws.Range("G10").Value2 = "1000 Pa"
$ws.Range("H10").Value2 = "kPa"
$ws.Range("G11").Value2 = "1000 Pascal"
$ws.Range("G12").Value2 = "*10^3 Pa"
$ws.Range("G13").Value2 = "1000000 Pa"
$ws.Range("G14").Value2 = "1000000 Pascal"
$ws.Range("G15").Value2 = "*10^6 Pa"
$ws.Range("H13").Value2 = "MPa"
$ws.Range("G16").Value2 = "1000000000 Pa"
$ws.Range("G17").Value2 = "1000000000 Pascal"
$ws.Range("G18").Value2 = "*10^9 Pa"
$ws.Range("H16").Value2 = "GPa"

# Step 2: fill in the remaining cells of the new lookup rows (2-18).
$ws.Range("M2").Value2 = "*10^3Hz"
$ws.Range("N2").Value2 = "kHz"
$ws.Range("O2").Value2 = "Frequency"
$ws.Range("H3").Value2 = "kHz"
$ws.Range("I3").Value2 = "Frequency"
$ws.Range("M3").Value2 = "1000 Hz"
$ws.Range("N3").Value2 = "kHz"
$ws.Range("O3").Value2 = "Frequency"
$ws.Range("I4").Value2 = "Frequency"
$ws.Range("N4").Value2 = "MHz"
$ws.Range("O4").Value2 = "Frequency"
$ws.Range("H5").Value2 = "MHz"
$ws.Range("I5").Value2 = "Frequency"
$ws.Range("M5").Value2 = "1000000 Hz"
$ws.Range("N5").Value2 = "MHz"
$ws.Range("O5").Value2 = "Frequency"
$ws.Range("I6").Value2 = "Frequency"
$ws.Range("M6").Value2 = "*10^9 Hz"
$ws.Range("N6").Value2 = "GHz"
$ws.Range("O6").Value2 = "Frequency"
$ws.Range("H7").Value2 = "GHz"
$ws.Range("I7").Value2 = "Frequency"
$ws.Range("M7").Value2 = "1000000000 Hz"
$ws.Range("N7").Value2 = "GHz"
$ws.Range("O7").Value2 = "Frequency"
$ws.Range("I8").Value2 = "Frequency"
$ws.Range("M8").Value2 = "rad/second"
$ws.Range("N8").Value2 = "rad*s^{-1}"
$ws.Range("O8").Value2 = "Frequency"
$ws.Range("M9").Value2 = "Pascal"
$ws.Range("N9").Value2 = "Pa"
$ws.Range("O9").Value2 = "Pressure"
$ws.Range("I10").Value2 = "Pressure"
$ws.Range("M10").Value2 = "1000 Pa"
$ws.Range("N10").Value2 = "kPa"
$ws.Range("O10").Value2 = "Pressure"
$ws.Range("H11").Value2 = "kPa"
$ws.Range("I11").Value2 = "Pressure"
$ws.Range("M11").Value2 = "1000 Pascal"
$ws.Range("N11").Value2 = "kPa"
$ws.Range("O11").Value2 = "Pressure"
$ws.Range("H12").Value2 = "kPa"
$ws.Range("I12").Value2 = "Pressure"
$ws.Range("M12").Value2 = "*10^3 Pa"
$ws.Range("N12").Value2 = "kPa"
$ws.Range("O12").Value2 = "Pressure"
$ws.Range("I13").Value2 = "Pressure"
$ws.Range("M13").Value2 = "1000000 Pa"
$ws.Range("N13").Value2 = "MPa"
$ws.Range("O13").Value2 = "Pressure"
$ws.Range("H14").Value2 = "MPa"
$ws.Range("I14").Value2 = "Pressure"
$ws.Range("M14").Value2 = "1000000 Pascal"
$ws.Range("N14").Value2 = "MPa"
$ws.Range("O14").Value2 = "Pressure"
$ws.Range("H15").Value2 = "MPa"
$ws.Range("I15").Value2 = "Pressure"
$ws.Range("M15").Value2 = "*10^6 Pa"
$ws.Range("N15").Value2 = "MPa"
$ws.Range("O15").Value2 = "Pressure"
$ws.Range("I16").Value2 = "Pressure"
$ws.Range("M16").Value2 = "1000000000 Pa"
$ws.Range("N16").Value2 = "GPa"
$ws.Range("O16").Value2 = "Pressure"
$ws.Range("H17").Value2 = "GPa"
$ws.Range("I17").Value2 = "Pressure"
$ws.Range("M17").Value2 = "1000000000 Pascal"
$ws.Range("N17").Value2 = "GPa"
$ws.Range("O17").Value2 = "Pressure"
$ws.Range("H18").Value2 = "GPa"
$ws.Range("I18").Value2 = "Pressure"
$ws.Range("M18").Value2 = "*10^9 Pa"
$ws.Range("N18").Value2 = "GPa"
$ws.Range("O18").Value2 = "Pressure"

# Update sheet view: scroll position and active selection in the frozen (bottom-left) pane
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("J10").Select()
